# Auto-generated edit script: update F-column ("想去人数") counts
# to match the refreshed scrape output (commit 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value2 = 1337
$ws.Range("F5").Value2 = 935
$ws.Range("F6").Value2 = 748
$ws.Range("F7").Value2 = 219
$ws.Range("F8").Value2 = 552
$ws.Range("F9").Value2 = 159
$ws.Range("F12").Value2 = 3096
$ws.Range("F13").Value2 = 2706
$ws.Range("F15").Value2 = 38
$ws.Range("F17").Value2 = 332
$ws.Range("F18").Value2 = 268
$ws.Range("F20").Value2 = 5537
$ws.Range("F21").Value2 = 605
$ws.Range("F24").Value2 = 72
$ws.Range("F25").Value2 = 421
$ws.Range("F26").Value2 = 1188
$ws.Range("F28").Value2 = 99
$ws.Range("F30").Value2 = 41

$ws = $wb.Worksheets.Item(2)
$ws.Range("F8").Value2 = 333
$ws.Range("F21").Value2 = 49
$ws.Range("F23").Value2 = 334
$ws.Range("F25").Value2 = 4023
$ws.Range("F26").Value2 = 3
$ws.Range("F29").Value2 = 206

$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value2 = 2540
$ws.Range("F6").Value2 = 1099
$ws.Range("F9").Value2 = 1409
$ws.Range("F10").Value2 = 396

$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value2 = 2540
$ws.Range("F6").Value2 = 1099
$ws.Range("F7").Value2 = 1409
$ws.Range("F8").Value2 = 396
$ws.Range("F11").Value2 = 1337
$ws.Range("F12").Value2 = 935
$ws.Range("F13").Value2 = 748
$ws.Range("F15").Value2 = 219
$ws.Range("F16").Value2 = 552
$ws.Range("F17").Value2 = 159
$ws.Range("F18").Value2 = 3096
$ws.Range("F19").Value2 = 2706
$ws.Range("F20").Value2 = 38
$ws.Range("F22").Value2 = 332
$ws.Range("F24").Value2 = 268
$ws.Range("F26").Value2 = 5537
$ws.Range("F28").Value2 = 605
$ws.Range("F32").Value2 = 72
$ws.Range("F33").Value2 = 421
$ws.Range("F38").Value2 = 49
$ws.Range("F39").Value2 = 334
$ws.Range("F40").Value2 = 1188
$ws.Range("F44").Value2 = 206
$ws.Range("F48").Value2 = 99
$ws.Range("F50").Value2 = 41
